$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44432
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 2300
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = 2400
$ws.Range("P2").Value = 2400

# Row 3
$ws.Range("D3").Value = 44365
$ws.Range("K3").Value = 2400
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = 2450
$ws.Range("P3").Value = 2450

# Row 4
$ws.Range("D4").Value = 44571
$ws.Range("K4").Value = 900
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = 950
$ws.Range("P4").Value = 950

# Row 5
$ws.Range("D5").Value = 44292
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 1900
$ws.Range("P5").Value = 1900

# Row 6
$ws.Range("D6").Value = 44539
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 900
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = 950
$ws.Range("P6").Value = 950

# Row 7
$ws.Range("D7").Value = 44498
$ws.Range("J7").Value = 270
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 2300
$ws.Range("M7").Value = 2150
$ws.Range("P7").Value = 2150

# Row 8
$ws.Range("D8").Value = 44224
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 750
$ws.Range("L8").Value = 800
$ws.Range("M8").Value = 775
$ws.Range("P8").Value = 775

# Row 9
$ws.Range("D9").Value = 44349

# Row 10
$ws.Range("D10").Value = 44326
$ws.Range("K10").Value = 2700
$ws.Range("L10").Value = 2800
$ws.Range("M10").Value = 2750
$ws.Range("P10").Value = 2750

# Row 11
$ws.Range("D11").Value = 44435
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 2300
$ws.Range("M11").Value = 2400
$ws.Range("P11").Value = 2400

# Row 12
$ws.Range("D12").Value = 44417
$ws.Range("K12").Value = 4000
$ws.Range("L12").Value = 4500
$ws.Range("M12").Value = 4250
$ws.Range("P12").Value = 4250

# Row 13
$ws.Range("D13").Value = 44362
$ws.Range("J13").Value = 250
$ws.Range("K13").Value = 2800
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = 2900
$ws.Range("P13").Value = 2900

# Row 14
$ws.Range("D14").Value = 44376
$ws.Range("J14").Value = 270
$ws.Range("K14").Value = 2400
$ws.Range("M14").Value = 2437
$ws.Range("P14").Value = 2437

# Row 15
$ws.Range("D15").Value = 44313
$ws.Range("J15").Value = 250

# Row 16
$ws.Range("D16").Value = 44540
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = 950
$ws.Range("P16").Value = 950

# Row 17
$ws.Range("D17").Value = 44280
$ws.Range("J17").Value = 250
$ws.Range("K17").Value = 1400
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = 1450
$ws.Range("P17").Value = 1450

# Row 18
$ws.Range("D18").Value = 44532
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 1000
$ws.Range("L18").Value = 1200
$ws.Range("M18").Value = 1100
$ws.Range("P18").Value = 1100

# Row 19
$ws.Range("D19").Value = 44442
$ws.Range("J19").Value = 240
$ws.Range("K19").Value = 2300
$ws.Range("L19").Value = 2500
$ws.Range("M19").Value = 2400
$ws.Range("P19").Value = 2400

# Row 20
$ws.Range("D20").Value = 44536
$ws.Range("K20").Value = 900
$ws.Range("L20").Value = 1000
$ws.Range("M20").Value = 950
$ws.Range("P20").Value = 950

# Row 21
$ws.Range("D21").Value = 44249
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 900
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 950
$ws.Range("P21").Value = 950

# Row 22
$ws.Range("D22").Value = 44330
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 2800
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = 2900
$ws.Range("P22").Value = 2900

# Row 23
$ws.Range("D23").Value = 44494
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 2400
$ws.Range("L23").Value = 2500
$ws.Range("M23").Value = 2450
$ws.Range("P23").Value = 2450

# Row 24
$ws.Range("D24").Value = 44250
$ws.Range("K24").Value = 1000
$ws.Range("L24").Value = 1200
$ws.Range("M24").Value = 1100
$ws.Range("P24").Value = 1100

# Row 25
$ws.Range("D25").Value = 44523
$ws.Range("K25").Value = 1400
$ws.Range("L25").Value = 1500
$ws.Range("M25").Value = 1450
$ws.Range("P25").Value = 1450

# Row 26
$ws.Range("D26").Value = 44260
$ws.Range("K26").Value = 900
$ws.Range("L26").Value = 1000
$ws.Range("M26").Value = 950
$ws.Range("P26").Value = 950

# Row 27
$ws.Range("D27").Value = 44474
$ws.Range("J27").Value = 250
$ws.Range("K27").Value = 2000
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = 2250
$ws.Range("P27").Value = 2250

# Row 28
$ws.Range("D28").Value = 44302
$ws.Range("K28").Value = 900
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 950
$ws.Range("P28").Value = 950

# Row 29
$ws.Range("D29").Value = 44274
$ws.Range("K29").Value = 1000
$ws.Range("L29").Value = 1200
$ws.Range("M29").Value = 1100
$ws.Range("P29").Value = 1100
